$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'59.395.05"
$ws.Range("E2").Value = "  +0.08%  "

# Row 3
$ws.Range("D3").Value = "'2.522.72"
$ws.Range("E3").Value = "  +0.14%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'535.80"
$ws.Range("E5").Value = "  -1.26%  "

# Row 6
$ws.Range("D6").Value = "'139.65"
$ws.Range("E6").Value = "  -3.60%  "

# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.34%  "

# Row 8
$ws.Range("D8").Value = "'0.564"
$ws.Range("E8").Value = "  -1.96%  "

# Row 9
$ws.Range("D9").Value = "'2.530.18"
$ws.Range("E9").Value = "  -0.62%  "

# Row 10
$ws.Range("E10").Value = "  +0.10%  "

# Row 11
$ws.Range("E11").Value = "  +1.25%  "

# Row 12
$ws.Range("D12").Value = "'5.47"
$ws.Range("E12").Value = "  -2.25%  "

# Row 13
$ws.Range("E13").Value = "  -0.01%  "

# Row 14
$ws.Range("D14").Value = "'2.971.21"
$ws.Range("E14").Value = "  +0.31%  "

# Row 15
$ws.Range("D15").Value = "'23.52"
$ws.Range("E15").Value = "  -1.42%  "

# Row 16
$ws.Range("D16").Value = "'59.306.13"
$ws.Range("E16").Value = "  +0.09%  "

# Row 17
$ws.Range("E17").Value = "  -0.04%  "

# Row 18
$ws.Range("D18").Value = "'2.522.33"
$ws.Range("E18").Value = "  -0.23%  "

# Row 19
$ws.Range("D19").Value = "'11.11"
$ws.Range("E19").Value = "  -1.48%  "

# Row 20
$ws.Range("D20").Value = "'4.32"
$ws.Range("E20").Value = "  +0.29%  "

# Row 21
$ws.Range("D21").Value = "'325.70"
$ws.Range("E21").Value = "  -0.35%  "

# Row 22
$ws.Range("E22").Value = "  +0.35%  "

# Row 23
$ws.Range("E23").Value = "  -0.33%  "

# Row 24
$ws.Range("D24").Value = "'63.73"
$ws.Range("E24").Value = "  +2.57%  "

# Row 25
$ws.Range("E25").Value = "  -2.16%  "

# Row 26
$ws.Range("E26").Value = "  +1.29%  "

# Row 27
$ws.Range("E27").Value = "  +0.90%  "

# Row 28
$ws.Range("E28").Value = "  -2.45%  "

# Row 29
$ws.Range("D29").Value = "'6.94"
$ws.Range("E29").Value = "  +1.13%  "

# Row 30
$ws.Range("D30").Value = "'0.0₃0780"
$ws.Range("E30").Value = "  -0.69%  "

# Row 31
$ws.Range("E31").Value = "  -2.65%  "

# Row 32
$ws.Range("D32").Value = "'165.43"
$ws.Range("E32").Value = "  +5.24%  "

# Row 33
$ws.Range("E33").Value = "  -1.74%  "

# Row 34
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  +0.15%  "

# Row 35
$ws.Range("E35").Value = "  -9.63%  "

# Row 36
$ws.Range("E36").Value = "  -0.80%  "

# Row 37
$ws.Range("E37").Value = "  -2.57%  "

# Row 38
$ws.Range("E38").Value = "  -1.86%  "

# Row 39
$ws.Range("D39").Value = "'36.96"
$ws.Range("E39").Value = "  +0.04%  "

# Row 40
$ws.Range("E40").Value = "  -0.59%  "

# Row 41
$ws.Range("E41").Value = "  -1.50%  "

# Row 42
$ws.Range("D42").Value = "'5.28"
$ws.Range("E42").Value = "  -6.92%  "

# Row 43
$ws.Range("D43").Value = "'280.59"
$ws.Range("E43").Value = "  -6.46%  "

# Row 44
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  +0.60%  "

# Row 45
$ws.Range("B45").Value = "'WhiteBITCoin"
$ws.Range("C45").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "'10.88"
$ws.Range("E45").Value = "  +0.75%  "

# Row 46
$ws.Range("B46").Value = "'Mantle"
$ws.Range("C46").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "'0.599"
$ws.Range("E46").Value = "  -0.97%  "

# Row 47
$ws.Range("E47").Value = "  -0.01%  "

# Row 48
$ws.Range("D48").Value = "'123.05"
$ws.Range("E48").Value = "  -0.49%  "

# Row 49
$ws.Range("E49").Value = "  -0.34%  "

# Row 50
$ws.Range("E50").Value = "  -1.67%  "

# Row 51
$ws.Range("D51").Value = "'17.84"
$ws.Range("E51").Value = "  -3.09%  "
